$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "285.38"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.45%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.66"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.60%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.101"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.91%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06677"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.69%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.333"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.21%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.383"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.48%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.354"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "4.49%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9361"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.95%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1574"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.72%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06502"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "14.23%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07641"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "2.68%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02885"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.81%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.08977"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.05%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001589"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.05%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04468"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.96%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006444"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.24%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006510"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "6.24%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.481"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.31%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.229"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.12%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3202"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.00%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.20%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.051"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "3.66%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.20%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001177"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.11%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004467"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.38%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.64%"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-2.40%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04186"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.01%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006729"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.12%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1245"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-10.59%"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.26%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01205"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "8.66%"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.58%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "20.74%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.51%"
